# Update 3rd order WENO coefficients calculator
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shift the right-hand "Order" block (rows 29:32) left by one column:
# old J,K,L,M -> new I,J,K,L
$ws.Range("I29:M32").Clear()

$ws.Cells.Item(29, 9).Value  = 2
$ws.Cells.Item(29, 10).Value = 1
$ws.Cells.Item(29, 11).Value = 0.5
$ws.Cells.Item(29, 12).Value = 0.25

$ws.Cells.Item(30, 9).Formula  = "=LOG(B30/C30)/LOG(B29/C29)"
$ws.Cells.Item(30, 10).Formula = "=LOG(C30/D30)/LOG(2)"
$ws.Cells.Item(30, 11).Formula = "=LOG(D30/E30)/LOG(2)"
$ws.Cells.Item(30, 12).Formula = "=LOG(E30/F30)/LOG(2)"
$ws.Range("I30:L30").NumberFormat = "0.00E+00"

$ws.Cells.Item(31, 9).Formula  = "=LOG(B31/C31)/LOG(B29/C29)"
$ws.Cells.Item(31, 10).Formula = "=LOG(C31/D31)/LOG(2)"
$ws.Cells.Item(31, 11).Formula = "=LOG(D31/E31)/LOG(2)"
$ws.Cells.Item(31, 12).Formula = "=LOG(E31/F31)/LOG(2)"
$ws.Range("I31:L31").NumberFormat = "0.00E+00"

$ws.Cells.Item(32, 9).Formula  = "=LOG(B32/C32)/LOG(B29/C29)"
$ws.Cells.Item(32, 10).Formula = "=LOG(C32/D32)/LOG(2)"
$ws.Cells.Item(32, 11).Formula = "=LOG(D32/E32)/LOG(2)"
$ws.Cells.Item(32, 12).Formula = "=LOG(E32/F32)/LOG(2)"
$ws.Range("I32:L32").NumberFormat = "0.00E+00"

# --- New "WENO" section (rows 34:38), mirrors the "Order" block above
$ws.Cells.Item(34, 1).Value2 = "WENO"

$ws.Cells.Item(35, 1).Value2 = "Order"
$ws.Cells.Item(35, 2).Value = 4.5
$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = 0.5
$ws.Cells.Item(35, 6).Value = 0.25
$ws.Cells.Item(35, 8).Value2 = "AUSM_2"
$ws.Cells.Item(35, 9).Value = 2
$ws.Cells.Item(35, 10).Value = 1
$ws.Cells.Item(35, 11).Value = 0.5
$ws.Cells.Item(35, 12).Value = 0.25

$ws.Cells.Item(36, 1).Value2 = "L1"
$ws.Cells.Item(36, 2).Value = 0.0000078192776413125699
$ws.Cells.Item(36, 3).Value = 0.000000108124282793362
$ws.Cells.Item(36, 4).Value = 0.0000000036421202685985199
$ws.Cells.Item(36, 5).Value = 0.000000000160302015969463
$ws.Cells.Item(36, 6).Value = 0.0000000000085089984231142698
$ws.Range("B36:F36").NumberFormat = "0.00E+00"
$ws.Cells.Item(36, 8).Value2 = "L1"
$ws.Cells.Item(36, 9).Formula  = "=LOG(B36/C36)/LOG(B35/C35)"
$ws.Cells.Item(36, 10).Formula = "=LOG(C36/D36)/LOG(2)"
$ws.Cells.Item(36, 11).Formula = "=LOG(D36/E36)/LOG(2)"
$ws.Cells.Item(36, 12).Formula = "=LOG(E36/F36)/LOG(2)"
$ws.Range("I36:L36").NumberFormat = "0.00E+00"

$ws.Cells.Item(37, 1).Value2 = "L2"
$ws.Cells.Item(37, 2).Value = 0.000010676656261578601
$ws.Cells.Item(37, 3).Value = 0.000000145635883906296
$ws.Cells.Item(37, 4).Value = 0.0000000046658453459182704
$ws.Cells.Item(37, 5).Value = 0.000000000190415033037408
$ws.Cells.Item(37, 6).Value = 0.0000000000100399432349302
$ws.Range("B37:F37").NumberFormat = "0.00E+00"
$ws.Cells.Item(37, 8).Value2 = "L2"
$ws.Cells.Item(37, 9).Formula  = "=LOG(B37/C37)/LOG(B35/C35)"
$ws.Cells.Item(37, 10).Formula = "=LOG(C37/D37)/LOG(2)"
$ws.Cells.Item(37, 11).Formula = "=LOG(D37/E37)/LOG(2)"
$ws.Cells.Item(37, 12).Formula = "=LOG(E37/F37)/LOG(2)"
$ws.Range("I37:L37").NumberFormat = "0.00E+00"

$ws.Cells.Item(38, 1).Value2 = "Linf"
$ws.Cells.Item(38, 2).Value = 0.000021363851494202999
$ws.Cells.Item(38, 3).Value = 0.00000027809675945248101
$ws.Cells.Item(38, 4).Value = 0.00000000809095568460738
$ws.Cells.Item(38, 5).Value = 0.00000000036509663150453699
$ws.Cells.Item(38, 6).Value = 0.000000000018631351385503401
$ws.Range("B38:F38").NumberFormat = "0.00E+00"
$ws.Cells.Item(38, 8).Value2 = "Linf"
$ws.Cells.Item(38, 9).Formula  = "=LOG(B38/C38)/LOG(B35/C35)"
$ws.Cells.Item(38, 10).Formula = "=LOG(C38/D38)/LOG(2)"
$ws.Cells.Item(38, 11).Formula = "=LOG(D38/E38)/LOG(2)"
$ws.Cells.Item(38, 12).Formula = "=LOG(E38/F38)/LOG(2)"
$ws.Range("I38:L38").NumberFormat = "0.00E+00"

# --- New "diff ratio" rows (40:42) comparing WENO (36:38) vs Order (30:32)
$ws.Cells.Item(40, 3).Formula = "=(C36-C30)/C30*100"
$ws.Cells.Item(40, 4).Formula = "=(D36-D30)/D30*100"
$ws.Range("C40:D40").NumberFormat = "0.00E+00"

$ws.Cells.Item(41, 3).Formula = "=(C37-C31)/C31*100"
$ws.Cells.Item(41, 4).Formula = "=(D37-D31)/D31*100"
$ws.Range("C41:D41").NumberFormat = "0.00E+00"

$ws.Cells.Item(42, 3).Formula = "=(C38-C32)/C32*100"
$ws.Cells.Item(42, 4).Formula = "=(D38-D32)/D32*100"
$ws.Range("C42:D42").NumberFormat = "0.00E+00"

# --- Column widths: col 2:4 used to be one merged block; now split so col 3 has its own width
$ws.Columns.Item(3).ColumnWidth = 9.285714285714286

# --- View: selection moves to J40, scrolled so row 22 is at the top
$ws.Range("J40").Select()
$excel.ActiveWindow.ScrollRow = 22

Write-Output "applied"
